# BOM PCB.xlsx edit: "export pcb1 aangepast, BOM PCB aangepast"
#
# The BOM table lost the separate "R2, R3 / 4k7 0805" line (merged into the
# "R1, R6, R7 / 10k 0805" line, which becomes "R1,R2,R3"), and the old
# "U7 / PA1010-D" + "L1 / MFBM1V1005-501-R" GPS-module lines were replaced by
# an "Adafruit Mini GPS" line (already bought - no link, just a note) and a
# "female pinheader 6" line.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the old "R2, R3 / 4k7 0805" row (row 10). Everything below it
#    (rows 11-17, plus the trailing formatted-but-empty row 21) shifts up by
#    one row.
# ---------------------------------------------------------------------------
$ws.Rows.Item(10).EntireRow.Delete()

# ---------------------------------------------------------------------------
# 2. Small text fixes on rows that stayed in place (values only, no
#    structural change).
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "U2,U6,"
$ws.Range("A6").Value = "C1,C2,C6,C8,C9,C11"
$ws.Range("A7").Value = "C4,C5,C7,C13"
$ws.Range("A9").Value = "R1,R2,R3"

# ---------------------------------------------------------------------------
# 3. Row 15 (was "U7 / PA1010-D", now shifted here): becomes the Adafruit
#    Mini GPS module, already bought -> no hyperlink, just a remark in G.
# ---------------------------------------------------------------------------
$ws.Range("B15").Value = "Adafruit Mini GPS"
$ws.Range("D15").Value = "mouser"
$ws.Range("E15").ClearContents()
$ws.Range("E15").Style = "Hyperlink"
$ws.Range("G15").Value = "al gekocht"

# ---------------------------------------------------------------------------
# 4. Row 16 (was "L1 / MFBM1V1005-501-R", now shifted here): becomes the
#    "female pinheader 6" connector row.
# ---------------------------------------------------------------------------
$ws.Range("A16").ClearContents()
$ws.Range("B16").Value = "female pinheader 6"
$ws.Range("C16").Value = 2
$ws.Range("F16").ClearContents()

# ---------------------------------------------------------------------------
# 5. Hyperlinks: row deletion does not re-anchor the existing hyperlink
#    objects to their new rows, so rebuild the whole collection from
#    scratch with the correct ref -> target pairing.
# ---------------------------------------------------------------------------
$ws.Range("E2").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("E2"), "https://www.mouser.be/ProductDetail/Espressif-Systems/ESP32-WROOM-32EM113EH2800PH3Q0?qs=sGAEpiMZZMu3sxpa5v1qrl%2FYtpu2q02IVMmJKBwXLzE%3D")
$ws.Hyperlinks.Add($ws.Range("E3"), "https://www.mouser.be/ProductDetail/LPRS/SMA-CONNECTOR?qs=j%252B1pi9TdxUYkOiITvzJM8A%3D%3D")
$ws.Hyperlinks.Add($ws.Range("E4"), "https://www.mouser.be/ProductDetail/Microchip-Technology/RN2483A-I-RM105?qs=sGAEpiMZZMu3sxpa5v1qrrZAnp88aoEUcCMjG6HT%2FSI%3D")
$ws.Hyperlinks.Add($ws.Range("E5"), "https://www.mouser.be/ProductDetail/Texas-Instruments/TPS22860DBVR?qs=%2Fha2pyFadugxAdFBsj4zSHcoO0BzEMtLRhA2IvKn%252BB5VwA7fNnvmmg%3D%3D")
$ws.Hyperlinks.Add($ws.Range("E6"), "https://www.mouser.be/ProductDetail/KEMET/C0805C104K9RACAUTO?qs=ds50AKTGxA84SWre%252BvImOA%3D%3D")
$ws.Hyperlinks.Add($ws.Range("E7"), "https://www.mouser.be/ProductDetail/KEMET/C0805X105K8RAC7210?qs=ds50AKTGxA8Ac4j4qbsB5A%3D%3D")
$ws.Hyperlinks.Add($ws.Range("E10"), "https://www.mouser.be/ProductDetail/Torex-Semiconductor/XC6231A332PR-G?qs=AsjdqWjXhJ%252BxcfXSnf052A%3D%3D")
$ws.Hyperlinks.Add($ws.Range("E11"), "https://www.mouser.be/ProductDetail/E-Switch/LL3301NF065QG?qs=sGAEpiMZZMsqIr59i2oRcpD0bdFgxAIu%2FtEwYcOnPMA%3D")
$ws.Hyperlinks.Add($ws.Range("E8"), "https://www.mouser.be/ProductDetail/Lelon/VE-220M1VTR-0605?qs=sGAEpiMZZMukHu%252BjC5l7YREIIWaVdOqUb1trgoTwpWE%3D")
$ws.Hyperlinks.Add($ws.Range("E12"), "https://www.mouser.be/ProductDetail/Molex/39700-0205?qs=sGAEpiMZZMvPvGwLNS6715dPPKD79FMV4ASqaQtPDLdJRa6HXWagoQ%3D%3D")
$ws.Hyperlinks.Add($ws.Range("E13"), "https://www.mouser.be/ProductDetail/Molex/39700-0002?qs=EibfsPFbZaqB1%252BjDW0HuPQ==")
$ws.Hyperlinks.Add($ws.Range("E14"), "https://www.mouser.be/ProductDetail/Gravitech/8Fx1L-254mm?qs=fkzBJ5HM%252BdAyuablm941Ag%3D%3D")
$ws.Hyperlinks.Add($ws.Range("E9"), "https://www.mouser.be/ProductDetail/Bourns/CMP0805AFX-1002ELF?qs=TiOZkKH1s2RdWJm6UffAxg%3D%3D")
$ws.Hyperlinks.Add($ws.Range("E16"), "https://www.mouser.be/ProductDetail/Bussmann-Eaton/MFBM1V1005-501-R?qs=%2Fha2pyFaduhkIz0cVTg91Fc3bMtUWIidpGzv3%2F1prZyfRm24uM5%252BKlDf5IXr9E00")

# Re-apply the Hyperlink cell style (Hyperlinks.Add() resets it) to every
# linked cell so the column keeps its original look.
foreach ($r in @(2,3,4,5,6,7,8,9,10,11,12,13,14,16)) {
    $ws.Range("E$r").Style = "Hyperlink"
}

# ---------------------------------------------------------------------------
# 6. Selection, matching the saved cursor position in the edited file.
# ---------------------------------------------------------------------------
$ws.Range("B17").Select()
